# LOT2039.xlsx edit — reshuffle the "Docentes/Programa/Syllabus/Avaliação" block
# (rows 10-22) per the target diff, then drop the now-unused trailing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Prime brand-new cells (A13, B18, C18) with the right look (bold label /
#        plain value / red value) by copying format from stable donor cells
#        whose rows are untouched by this edit (row 3). Copy() brings the cell
#        style along; the actual text is overwritten right after.
$ws.Range("A3").Copy($ws.Range("A13"))
$ws.Range("B3").Copy($ws.Range("B18"))
$ws.Range("C3").Copy($ws.Range("C18"))

# --- 2) Content-preserving moves. These must run BEFORE the cells that serve
#        as their source get overwritten later in the script.

# B13/C13 ("3380737 - Flávio Teixeira da Silva") -> B10/C10 and B18/C18
$ws.Range("B13").Copy($ws.Range("B10"))
$ws.Range("C13").Copy($ws.Range("C10"))
$ws.Range("B13").Copy($ws.Range("B18"))
$ws.Range("C13").Copy($ws.Range("C18"))

# B15/C15 (Short syllabus EN text) -> B14/C14
$ws.Range("B15").Copy($ws.Range("B14"))
$ws.Range("C15").Copy($ws.Range("C14"))

# B17/C17 (full Syllabus EN text) -> B16/C16
$ws.Range("B17").Copy($ws.Range("B16"))
$ws.Range("C17").Copy($ws.Range("C16"))

# B8/C8 ("01/01/2018") -> B15/C15 (B8/C8 itself is never modified, so order-safe)
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# --- 3) Brand-new literal text.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("A15").Value = "Programa:"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A21").Value = "Bibliografia:"

# --- 4) Drop the old long "Programa:" text — row 17 becomes an A-only row.
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()

# --- 5) Row heights that changed.
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120

# --- 6) The table lost a row overall; drop the trailing (now orphaned) row 22.
$ws.Rows.Item(22).Delete()
